$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-7 (A..T), reflecting the new TPM values and the
# addition of a "Resolving-Mac" target cluster alongside the existing
# "MuSCs" target cluster for each sending cluster (ECs, FAPs, MuSCs).
$data = @(
    @("ECs",  "Efna5", "Epha5", "MuSCs",         2, 0.6666666666666666, 0.3227736666666667, 0.968321,   0.1416094457286952,  0.1416094457286952,  3, 1,                  0.02440666666666666,  0.07321999999999999, 0.9884709884709885,  0.9884709884709886,  0.00787782929111111,    0.07090046361999999, 0.1399768287962721,    0.1399768287962721),
    @("ECs",  "Efna5", "Epha5", "Resolving-Mac", 2, 0.6666666666666666, 0.3227736666666667, 0.968321,   0.1416094457286952,  0.1416094457286952,  1, 0.3333333333333333, 0.0002846666666666667, 0.000854,             0.01152901152901153, 0.01152901152901153, 0.00009188290377777778, 0.0008269461340000001, 0.001632616932423059, 0.001632616932423059),
    @("FAPs", "Efna5", "Epha5", "MuSCs",         3, 1,                  1.874986333333333,  5.624959,   0.8226066833587575,  0.8226066833587576,  3, 1,                  0.02440666666666666,  0.07321999999999999, 0.9884709884709885,  0.9884709884709886,  0.04576216644222222,    0.41185949798,        0.8131228414224725,   0.8131228414224727),
    @("FAPs", "Efna5", "Epha5", "Resolving-Mac", 3, 1,                  1.874986333333333,  5.624959,   0.8226066833587575,  0.8226066833587576,  1, 0.3333333333333333, 0.0002846666666666667, 0.000854,             0.01152901152901153, 0.01152901152901153, 0.0005337461095555556,  0.004803714986000001, 0.009483841936285054, 0.009483841936285054),
    @("MuSCs","Efna5", "Epha5", "MuSCs",         2, 0.6666666666666666, 0.081563,            0.244689,   0.03578387091254728, 0.03578387091254728, 3, 1,                  0.02440666666666666,  0.07321999999999999, 0.9884709884709885,  0.9884709884709886,  0.001990680953333333,   0.01791612858,        0.03537131825224386,  0.03537131825224386),
    @("MuSCs","Efna5", "Epha5", "Resolving-Mac", 2, 0.6666666666666666, 0.081563,            0.244689,   0.03578387091254728, 0.03578387091254728, 1, 0.3333333333333333, 0.0002846666666666667, 0.000854,             0.01152901152901153, 0.01152901152901153, 0.00002321826733333333, 0.000208964406,       0.0004125526603034179, 0.0004125526603034179)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowValues[$j]
    }
}
